$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44159
$ws.Range("K2").Value = "Castle Brite"
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 14500
$ws.Range("Q2").Value = "$/caja 15 kilos"
$ws.Range("S2").Value = 967
$ws.Range("T2").Value = 15

# Row 3
$ws.Range("D3").Value = 44189
$ws.Range("K3").Value = "Dina"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 15500
$ws.Range("Q3").Value = "$/caja 15 kilos granel"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1033

# Row 4
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 14000
$ws.Range("S4").Value = 933

# Row 5
$ws.Range("D5").Value = 44187
$ws.Range("L5").Value = "Primera"
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 15500
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 861
$ws.Range("T5").Value = 18
